# "Update van de tijdsbesteding: afgelopen 2 weken"
# Fill in the two still-open weeks (rows 6-8, column B) with this period's
# logged hours, which ripples into the totaal/gemiddeld formulas in E2/E3,
# and leave the selection on the newly-edited range (B6:B7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("week 16-22/11/2015"): previously just "=7", now several entries summed.
$ws.Range("B6").Formula = "=7+1+1+0.5+3"

# Row 7 ("week 09-15/11/2015"): was empty, now filled in.
$ws.Range("B7").Formula = "=4.75+3.33+2.5"

# Row 8 ("week 02-08/11/2015"): was empty, now filled in.
$ws.Range("B8").Formula = "=4"

# Reflect the Mac-Excel selection state captured after the edit: active
# cell B6, with B6:B7 highlighted.
$ws.Range("B6:B7").Select()
